# Helper: force a cell to hold a genuine *text* value, even when the
# string looks numeric (e.g. "4811", "1176", "3e00" would otherwise be
# auto-converted to numbers/scientific notation by plain .Value / .Value2
# assignment). Entering it as a quoted-text formula and then pasting back
# as values-only keeps the cell's existing (lack of) style untouched.
function Set-TextValue($range, [string]$value) {
    $range.Formula = "=""" + $value + """"
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

# --- Sheet1: refresh the randomized User ID codes ---
$ws1 = $wb.Worksheets.Item("Sheet1")
Set-TextValue $ws1.Range("A2") "2f3c"
Set-TextValue $ws1.Range("A3") "4811"
Set-TextValue $ws1.Range("A4") "c046"
Set-TextValue $ws1.Range("A5") "1176"

# --- Login sheet: same User ID codes, kept in sync ---
$ws2 = $wb.Worksheets.Item("Login")
Set-TextValue $ws2.Range("A2") "2f3c"
Set-TextValue $ws2.Range("A3") "4811"
Set-TextValue $ws2.Range("A4") "c046"
Set-TextValue $ws2.Range("A5") "1176"

# --- Products details: refresh the randomized Product ID codes ---
$ws3 = $wb.Worksheets.Item("Products details")
Set-TextValue $ws3.Range("A2") "3e00"
Set-TextValue $ws3.Range("A3") "08b2"
Set-TextValue $ws3.Range("A4") "568f"
Set-TextValue $ws3.Range("A5") "cc8b"
Set-TextValue $ws3.Range("A6") "d499"
Set-TextValue $ws3.Range("A7") "7959"

$excel.CutCopyMode = $false

# --- Order Details: add the new "Order Status" column (task 5 - place orders) ---
$ws4 = $wb.Worksheets.Item("Order Details")
$ws4.Range("G1").Value = "Order Status"
$ws4.Range("G2").Value = "Success"
$ws4.Range("G3").Value = "Success"
$ws4.Range("G4").Value = "Success"
